$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Set3")
$ws.Range("B2").Value = 9.0206898844950008
